$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.295.59"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.899.15"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - was BNB, now XRP
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "0.693"
$ws.Range("E5").Value = "  +9.44%  "

# Row 6 - was XRP, now BNB
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "246.27"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.13%  "

# Row 8 - Solana
$ws.Range("D8").Value = "40.39"
$ws.Range("E8").Value = "  -3.53%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.350"
$ws.Range("E9").Value = "  +3.49%  "

# Row 10 - OKB
$ws.Range("D10").Value = "51.96"
$ws.Range("E10").Value = "  +6.69%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.45%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.0987"
$ws.Range("E12").Value = "  -0.94%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.174.57"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "12.52"
$ws.Range("E14").Value = "  +1.54%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +2.65%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.899.32"
$ws.Range("E16").Value = "  -0.62%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "4.84"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "35.286.25"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").Value = "72.50"
$ws.Range("E19").Value = "  +0.69%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0819"
$ws.Range("E20").Value = "  -0.82%  "

# Row 21
$ws.Range("D21").Value = "241.02"
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("D22").Value = "12.80"
$ws.Range("E22").Value = "  +2.34%  "

# Row 23
$ws.Range("D23").Value = "5.03"
$ws.Range("E23").Value = "  +3.22%  "

# Row 24
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  +0.42%  "

# Row 26
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  +5.11%  "

# Row 27
$ws.Range("D27").Value = "167.69"
$ws.Range("E27").Value = "  -2.69%  "

# Row 28
$ws.Range("E28").Value = "  +0.16%  "

# Row 29
$ws.Range("D29").Value = "18.95"
$ws.Range("E29").Value = "  +5.43%  "

# Row 30
$ws.Range("D30").Value = "0.131"
$ws.Range("E30").Value = "  +4.46%  "

# Row 32
$ws.Range("E32").Value = "  +2.37%  "

# Row 33
$ws.Range("D33").Value = "0.0570"
$ws.Range("E33").Value = "  +0.09%  "

# Row 34
$ws.Range("E34").Value = "  -0.18%  "

# Row 35
$ws.Range("E35").Value = "  +6.85%  "

# Row 36
$ws.Range("E36").Value = "  -1.12%  "

# Row 37
$ws.Range("D37").Value = "0.914"
$ws.Range("E37").Value = "  -5.09%  "

# Row 38
$ws.Range("E38").Value = "  +7.98%  "

# Row 39
$ws.Range("E39").Value = "  -0.50%  "

# Row 40
$ws.Range("D40").Value = "95.89"
$ws.Range("E40").Value = "  +5.42%  "

# Row 41 - was Kaspa, now ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.09"
$ws.Range("E41").Value = "  -1.00%  "

# Row 42 - was ARBITRUM, now Kaspa
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.0652"
$ws.Range("E42").Value = "  +7.24%  "

# Row 43
$ws.Range("D43").Value = "16.37"
$ws.Range("E43").Value = "  +3.83%  "

# Row 44
$ws.Range("D44").Value = "0.0208"
$ws.Range("E44").Value = "  +1.20%  "

# Row 45
$ws.Range("D45").Value = "1.355.32"
$ws.Range("E45").Value = "  +0.25%  "

# Row 46
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("E48").Value = "  +0.78%  "

# Row 49
$ws.Range("D49").Value = "45.32"
$ws.Range("E49").Value = "  -8.34%  "

# Row 50
$ws.Range("D50").Value = "12.32"
$ws.Range("E50").Value = "  -3.15%  "

# Row 51
$ws.Range("D51").Value = "6.49"
$ws.Range("E51").Value = "  -2.65%  "
